# Delete specific rows from the "Export" sheet by account number (column A).
# The rows being removed correspond to:
#   004213929 - Rodolfo
#   008364902 - Marcio
#   000806386 - Fernanda
#   005103059 - Walquiria
#   005135532 - Felipe

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$accountsToDelete = @("004213929", "008364902", "000806386", "005103059", "005135532")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Collect the row numbers that match the target account numbers first,
# so we can delete them from bottom to top (avoids row-index shifting issues).
$rowsToDelete = New-Object System.Collections.ArrayList

for ($r = 1; $r -le $lastRow; $r++) {
    $cellValue = $ws.Cells.Item($r, 1).Text
    if ($accountsToDelete -contains $cellValue) {
        [void]$rowsToDelete.Add($r)
    }
}

$sorted = $rowsToDelete | Sort-Object -Descending

foreach ($r in $sorted) {
    $ws.Rows.Item($r).Delete()
}
